$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("M", 20),
    @("N", 22),
    @("O", 30),
    @("P", 32),
    @("Q", 45),
    @("R", 40),
    @("S", 44),
    @("T", 14),
    @("U", 38),
    @("V", 38),
    @("X", 11),
    @("W", 31),
    @("Y", 50),
    @("Z", 34)
)

$row = 7
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}

[void]$ws.Range("B7:B20").Select()
